$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.93632866666666
$ws.Range("H2").Value = 53.808986
$ws.Range("I2").Value = 0.1226979812530711
$ws.Range("J2").Value = 0.1347750935001359
$ws.Range("M2").Value = 2.621797333333333
$ws.Range("N2").Value = 7.865392
$ws.Range("O2").Value = 0.07867217155043885
$ws.Range("P2").Value = 0.07906089226781998
$ws.Range("Q2").Value = 47.02541866805689
$ws.Range("R2").Value = 423.228768012512
$ws.Range("S2").Value = 0.009652916630034142
$ws.Range("T2").Value = 0.01065543914759961
$ws.Range("G3").Value = 17.93632866666666
$ws.Range("H3").Value = 53.808986
$ws.Range("I3").Value = 0.1226979812530711
$ws.Range("J3").Value = 0.1347750935001359
$ws.Range("O3").Value = 0.148308476825081
$ws.Range("P3").Value = 0.1490412718702539
$ws.Range("Q3").Value = 88.64974840881132
$ws.Range("R3").Value = 797.847735679302
$ws.Range("S3").Value = 0.01819715070915533
$ws.Range("T3").Value = 0.02008705135169264
$ws.Range("G4").Value = 17.93632866666666
$ws.Range("H4").Value = 53.808986
$ws.Range("I4").Value = 0.1226979812530711
$ws.Range("J4").Value = 0.1347750935001359
$ws.Range("M4").Value = 11.09754033333333
$ws.Range("N4").Value = 33.292621
$ws.Range("O4").Value = 0.3330034651388949
$ws.Range("P4").Value = 0.3346488416844782
$ws.Range("Q4").Value = 199.0491308102562
$ws.Range("R4").Value = 1791.442177292306
$ws.Range("S4").Value = 0.04085885292281986
$ws.Range("T4").Value = 0.04510232892773771
$ws.Range("G5").Value = 17.93632866666666
$ws.Range("H5").Value = 53.808986
$ws.Range("I5").Value = 0.1226979812530711
$ws.Range("J5").Value = 0.1347750935001359
$ws.Range("M5").Value = 0.4915585
$ws.Range("N5").Value = 0.983117
$ws.Range("O5").Value = 0.01475017696730553
$ws.Range("P5").Value = 0.009882038584175128
$ws.Range("Q5").Value = 8.816754814893665
$ws.Range("R5").Value = 52.900528889362
$ws.Range("S5").Value = 0.001809816937013936
$ws.Range("T5").Value = 0.001331852674154153
$ws.Range("G6").Value = 17.93632866666666
$ws.Range("H6").Value = 53.808986
$ws.Range("I6").Value = 0.1226979812530711
$ws.Range("J6").Value = 0.1347750935001359
$ws.Range("M6").Value = 14.172235
$ws.Range("N6").Value = 42.516705
$ws.Range("O6").Value = 0.4252657095182797
$ws.Range("P6").Value = 0.4273669555932728
$ws.Range("Q6").Value = 254.1978649012366
$ws.Range("R6").Value = 2287.78078411113
$ws.Range("S6").Value = 0.05217924405404788
$ws.Range("T6").Value = 0.05759842139895176
$ws.Range("I7").Value = 0.3929554311523962
$ws.Range("J7").Value = 0.4316338739568692
$ws.Range("M7").Value = 2.621797333333333
$ws.Range("N7").Value = 7.865392
$ws.Range("O7").Value = 0.07867217155043885
$ws.Range("P7").Value = 0.07906089226781998
$ws.Range("Q7").Value = 150.6047082365156
$ws.Range("R7").Value = 1355.44237412864
$ws.Range("S7").Value = 0.03091465709129798
$ws.Range("T7").Value = 0.03412535920804582
$ws.Range("I8").Value = 0.3929554311523962
$ws.Range("J8").Value = 0.4316338739568692
$ws.Range("O8").Value = 0.148308476825081
$ws.Range("P8").Value = 0.1490412718702539
$ws.Range("S8").Value = 0.05827862145435488
$ws.Range("T8").Value = 0.06433126155681665
$ws.Range("I9").Value = 0.3929554311523962
$ws.Range("J9").Value = 0.4316338739568692
$ws.Range("M9").Value = 11.09754033333333
$ws.Range("N9").Value = 33.292621
$ws.Range("O9").Value = 0.3330034651388949
$ws.Range("P9").Value = 0.3346488416844782
$ws.Range("Q9").Value = 637.4794126133689
$ws.Range("R9").Value = 5737.31471352032
$ws.Range("S9").Value = 0.1308555202188964
$ws.Range("T9").Value = 0.1444457759514503
$ws.Range("I10").Value = 0.3929554311523962
$ws.Range("J10").Value = 0.4316338739568692
$ws.Range("M10").Value = 0.4915585
$ws.Range("N10").Value = 0.983117
$ws.Range("O10").Value = 0.01475017696730553
$ws.Range("P10").Value = 0.009882038584175128
$ws.Range("Q10").Value = 28.23674566010667
$ws.Range("R10").Value = 169.42047396064
$ws.Range("S10").Value = 0.005796162149761689
$ws.Range("T10").Value = 0.004265422596678765
$ws.Range("I11").Value = 0.3929554311523962
$ws.Range("J11").Value = 0.4316338739568692
$ws.Range("M11").Value = 14.172235
$ws.Range("N11").Value = 42.516705
$ws.Range("O11").Value = 0.4252657095182797
$ws.Range("P11").Value = 0.4273669555932728
$ws.Range("Q11").Value = 814.1000412570668
$ws.Range("R11").Value = 7326.900371313601
$ws.Range("S11").Value = 0.1671104702380853
$ws.Range("T11").Value = 0.1844660546438776
$ws.Range("G12").Value = 9.626273333333334
$ws.Range("H12").Value = 28.87882
$ws.Range("I12").Value = 0.06585095126993876
$ws.Range("J12").Value = 0.07233263354328205
$ws.Range("M12").Value = 2.621797333333333
$ws.Range("N12").Value = 7.865392
$ws.Range("O12").Value = 0.07867217155043885
$ws.Range("P12").Value = 0.07906089226781998
$ws.Range("Q12").Value = 25.23813775527111
$ws.Range("R12").Value = 227.14323979744
$ws.Range("S12").Value = 0.00518063733506821
$ws.Range("T12").Value = 0.005718682548013125
$ws.Range("G13").Value = 9.626273333333334
$ws.Range("H13").Value = 28.87882
$ws.Range("I13").Value = 0.06585095126993876
$ws.Range("J13").Value = 0.07233263354328205
$ws.Range("O13").Value = 0.148308476825081
$ws.Range("P13").Value = 0.1490412718702539
$ws.Range("Q13").Value = 47.57755753552667
$ws.Range("R13").Value = 428.19801781974
$ws.Range("S13").Value = 0.009766254280327251
$ws.Range("T13").Value = 0.01078054770101575
$ws.Range("G14").Value = 9.626273333333334
$ws.Range("H14").Value = 28.87882
$ws.Range("I14").Value = 0.06585095126993876
$ws.Range("J14").Value = 0.07233263354328205
$ws.Range("M14").Value = 11.09754033333333
$ws.Range("N14").Value = 33.292621
$ws.Range("O14").Value = 0.3330034651388949
$ws.Range("P14").Value = 0.3346488416844782
$ws.Range("Q14").Value = 106.8279565763578
$ws.Range("R14").Value = 961.4516091872199
$ws.Range("S14").Value = 0.02192859495558212
$ws.Range("T14").Value = 0.02420603203124717
$ws.Range("G15").Value = 9.626273333333334
$ws.Range("H15").Value = 28.87882
$ws.Range("I15").Value = 0.06585095126993876
$ws.Range("J15").Value = 0.07233263354328205
$ws.Range("M15").Value = 0.4915585
$ws.Range("N15").Value = 0.983117
$ws.Range("O15").Value = 0.01475017696730553
$ws.Range("P15").Value = 0.009882038584175128
$ws.Range("Q15").Value = 4.731876480323334
$ws.Range("R15").Value = 28.39125888194
$ws.Range("S15").Value = 0.0009713131846970094
$ws.Range("T15").Value = 0.0007147938755697133
$ws.Range("G16").Value = 9.626273333333334
$ws.Range("H16").Value = 28.87882
$ws.Range("I16").Value = 0.06585095126993876
$ws.Range("J16").Value = 0.07233263354328205
$ws.Range("M16").Value = 14.172235
$ws.Range("N16").Value = 42.516705
$ws.Range("O16").Value = 0.4252657095182797
$ws.Range("P16").Value = 0.4273669555932728
$ws.Range("Q16").Value = 136.4258078542333
$ws.Range("R16").Value = 1227.8322706881
$ws.Range("S16").Value = 0.02800415151426417
$ws.Range("T16").Value = 0.03091257738743629
$ws.Range("G17").Value = 39.29803649999999
$ws.Range("H17").Value = 78.59607299999999
$ws.Range("I17").Value = 0.2688281328564436
$ws.Range("J17").Value = 0.1968591842135532
$ws.Range("M17").Value = 2.621797333333333
$ws.Range("N17").Value = 7.865392
$ws.Range("O17").Value = 0.07867217155043885
$ws.Range("P17").Value = 0.07906089226781998
$ws.Range("Q17").Value = 103.031487300936
$ws.Range("R17").Value = 618.1889238056159
$ws.Range("S17").Value = 0.0211492929856663
$ws.Range("T17").Value = 0.01556386275503865
$ws.Range("G18").Value = 39.29803649999999
$ws.Range("H18").Value = 78.59607299999999
$ws.Range("I18").Value = 0.2688281328564436
$ws.Range("J18").Value = 0.1968591842135532
$ws.Range("O18").Value = 0.148308476825081
$ws.Range("P18").Value = 0.1490412718702539
$ws.Range("Q18").Value = 194.2293271621185
$ws.Range("R18").Value = 1165.375962972711
$ws.Range("S18").Value = 0.03986949091166967
$ws.Range("T18").Value = 0.02934014319452857
$ws.Range("G19").Value = 39.29803649999999
$ws.Range("H19").Value = 78.59607299999999
$ws.Range("I19").Value = 0.2688281328564436
$ws.Range("J19").Value = 0.1968591842135532
$ws.Range("M19").Value = 11.09754033333333
$ws.Range("N19").Value = 33.292621
$ws.Range("O19").Value = 0.3330034651388949
$ws.Range("P19").Value = 0.3346488416844782
$ws.Range("Q19").Value = 436.1115450795554
$ws.Range("R19").Value = 2616.669270477332
$ws.Range("S19").Value = 0.08952069976801494
$ws.Range("T19").Value = 0.06587869797201688
$ws.Range("G20").Value = 39.29803649999999
$ws.Range("H20").Value = 78.59607299999999
$ws.Range("I20").Value = 0.2688281328564436
$ws.Range("J20").Value = 0.1968591842135532
$ws.Range("M20").Value = 0.4915585
$ws.Range("N20").Value = 0.983117
$ws.Range("O20").Value = 0.01475017696730553
$ws.Range("P20").Value = 0.009882038584175128
$ws.Range("Q20").Value = 19.31728387488525
$ws.Range("R20").Value = 77.26913549954099
$ws.Range("S20").Value = 0.003965262533422866
$ws.Range("T20").Value = 0.001945370054047572
$ws.Range("G21").Value = 39.29803649999999
$ws.Range("H21").Value = 78.59607299999999
$ws.Range("I21").Value = 0.2688281328564436
$ws.Range("J21").Value = 0.1968591842135532
$ws.Range("M21").Value = 14.172235
$ws.Range("N21").Value = 42.516705
$ws.Range("O21").Value = 0.4252657095182797
$ws.Range("P21").Value = 0.4273669555932728
$ws.Range("Q21").Value = 556.9410083165775
$ws.Range("R21").Value = 3341.646049899465
$ws.Range("S21").Value = 0.1143233866576699
$ws.Range("T21").Value = 0.08413111023792148
$ws.Range("G22").Value = 21.87880766666666
$ws.Range("H22").Value = 65.63642299999999
$ws.Range("I22").Value = 0.1496675034681502
$ws.Range("J22").Value = 0.1643992147861598
$ws.Range("M22").Value = 2.621797333333333
$ws.Range("N22").Value = 7.865392
$ws.Range("O22").Value = 0.07867217155043885
$ws.Range("P22").Value = 0.07906089226781998
$ws.Range("Q22").Value = 57.36179959697954
$ws.Range("R22").Value = 516.2561963728159
$ws.Range("S22").Value = 0.01177466750837221
$ws.Range("T22").Value = 0.01299754860912278
$ws.Range("G23").Value = 21.87880766666666
$ws.Range("H23").Value = 65.63642299999999
$ws.Range("I23").Value = 0.1496675034681502
$ws.Range("J23").Value = 0.1643992147861598
$ws.Range("O23").Value = 0.148308476825081
$ws.Range("P23").Value = 0.1490412718702539
$ws.Range("Q23").Value = 108.1353286494623
$ws.Range("R23").Value = 973.217957845161
$ws.Range("S23").Value = 0.02219695946957389
$ws.Range("T23").Value = 0.02450226806620031
$ws.Range("G24").Value = 21.87880766666666
$ws.Range("H24").Value = 65.63642299999999
$ws.Range("I24").Value = 0.1496675034681502
$ws.Range("J24").Value = 0.1643992147861598
$ws.Range("M24").Value = 11.09754033333333
$ws.Range("N24").Value = 33.292621
$ws.Range("O24").Value = 0.3330034651388949
$ws.Range("P24").Value = 0.3346488416844782
$ws.Range("Q24").Value = 242.8009505260758
$ws.Range("R24").Value = 2185.208554734682
$ws.Range("S24").Value = 0.04983979727358159
$ws.Range("T24").Value = 0.05501600680202613
$ws.Range("G25").Value = 21.87880766666666
$ws.Range("H25").Value = 65.63642299999999
$ws.Range("I25").Value = 0.1496675034681502
$ws.Range("J25").Value = 0.1643992147861598
$ws.Range("M25").Value = 0.4915585
$ws.Range("N25").Value = 0.983117
$ws.Range("O25").Value = 0.01475017696730553
$ws.Range("P25").Value = 0.009882038584175128
$ws.Range("Q25").Value = 10.75471387841517
$ws.Range("R25").Value = 64.52828327049099
$ws.Range("S25").Value = 0.00220762216241003
$ws.Range("T25").Value = 0.001624599383724926
$ws.Range("G26").Value = 21.87880766666666
$ws.Range("H26").Value = 65.63642299999999
$ws.Range("I26").Value = 0.1496675034681502
$ws.Range("J26").Value = 0.1643992147861598
$ws.Range("M26").Value = 14.172235
$ws.Range("N26").Value = 42.516705
$ws.Range("O26").Value = 0.4252657095182797
$ws.Range("P26").Value = 0.4273669555932728
$ws.Range("Q26").Value = 310.0716037718016
$ws.Range("R26").Value = 2790.644433946215
$ws.Range("S26").Value = 0.06364845705421249
$ws.Range("T26").Value = 0.07025879192508569
